$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

$ws.Range("B2").Value = 700
$ws.Range("B3").Value = 500
$ws.Range("B4").Value = 150
$ws.Range("B5").Value = 40
$ws.Range("B6").Value = 285
